$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Insert two new rows before current row 5 (pushing "Rebecca McLean" row and below down by 2)
$ws1.Rows("5:6").Insert()

# Fill in the newly inserted rows
$ws1.Cells.Item(5, 1).Value = 111
$ws1.Cells.Item(5, 2).Value = "Nod"
$ws1.Cells.Item(5, 3).Value = "Rod"
$ws1.Cells.Item(5, 4).Value = "Brazil"

$ws1.Cells.Item(6, 1).Value = 112
$ws1.Cells.Item(6, 2).Value = "Mori"
$ws1.Cells.Item(6, 3).Value = "Cox"
$ws1.Cells.Item(6, 4).Value = "Laos"

# Hide the two newly inserted rows
$ws1.Rows("5:6").Hidden = $true

# Sheet2 keeps its previous selection, but is no longer the active tab
$ws2.Range("J13").Select()

# Sheet1 becomes the active sheet/tab, with a new selection
$ws1.Activate()
$ws1.Range("E19").Select()
